$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values for Hong Kong Insurance (Prop/Cas.) rows 2 and 3

$updates = @{
    "D"  = 0.0582
    "G"  = -0.3818827708703375
    "H"  = -0.3818827708703375
    "I"  = -0.3552397868561279
    "J"  = -0.3552397868561279
    "K"  = -21.7
    "L"  = -0.3854351687388988
    "U"  = 110.6
    "V"  = 2.835897435897436
    "W"  = -0.3640939597315436
    "X"  = 0.06046991799701271
    "Y"  = -0.4245638777285563
    "Z"  = 2.336099585062241
    "AA" = -0.8298755186721993
    "AB" = 0.04800324368620691
    "AC" = -0.8778787623584062
    "AD" = 22.1
    "AF" = 22.1
    "AG" = -88.5
    "AH" = 0.3617021276595745
    "AI" = 0.3708053691275168
    "AJ" = 1.787878787878788
    "AK" = 1.735294117647059
    "AL" = 0.33
    "AM" = 0.33
    "AN" = -1.284883720930233
    "AO" = -60.6060606060606
    "AP" = 5.145348837209302
    "AQ" = -60.6060606060606
}

foreach ($row in @(2, 3)) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
